$d = $word.ActiveDocument

$old1 = "새로운 요리를 만들어보기 위해 혼자 재료를 준비하고 레시피를 찾아다니고, 피드백을 받을 수 없어서 만들어진 요리가 잘 된 건지 아닌 지도 모르게 되는 "
$new1 = "새로운 요리를 만들어보기 위해 혼자 막막하게 레시피를 찾아다니고, 없는 재료를 준비하고, 다 만들고 나서도 피드백 해줄 사람이 없어서 만들어진 요리가 잘 된 건지 아닌 지도 모르는 상황은 좀 많이 답답하기 때문에 올해 여름부터 꾸준히 요리학원에 등록해서 한식과 떡, 한과 디저트 제작 등을 배우고 있다. "

$d.Content.Find.Execute($old1, $true, $false, $false, $false, $false,
                         $true, 1, $false, $new1, 2)
